$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 500  # H5
$ws.Cells.Item(5, 9).Value = 750  # I5
$ws.Cells.Item(5, 11).Value = 750  # K5
$ws.Cells.Item(5, 13).Value = -635  # M5
$ws.Cells.Item(18, 8).Value = 15624.3125  # H18
$ws.Cells.Item(18, 9).Value = 14833.333  # I18
$ws.Cells.Item(18, 10).Value = 16098.9  # J18
$ws.Cells.Item(18, 11).Value = 14833.333  # K18
$ws.Cells.Item(18, 12).Value = 16098.9  # L18
$ws.Cells.Item(18, 13).Value = -14549.333  # M18
$ws.Cells.Item(18, 14).Value = -16666.9  # N18
$ws.Cells.Item(70, 8).Value = 90000  # H70
$ws.Cells.Item(70, 9).Value = 80000  # I70
$ws.Cells.Item(70, 10).Value = 100000  # J70
$ws.Cells.Item(70, 11).Value = 240000  # K70
$ws.Cells.Item(70, 12).Value = 300000  # L70
$ws.Cells.Item(70, 13).Value = -239730  # M70
$ws.Cells.Item(70, 14).Value = -300540  # N70
$ws.Cells.Item(73, 8).Value = 90000  # H73
$ws.Cells.Item(73, 9).Value = 80000  # I73
$ws.Cells.Item(73, 10).Value = 100000  # J73
$ws.Cells.Item(73, 11).Value = 240000  # K73
$ws.Cells.Item(73, 12).Value = 300000  # L73
$ws.Cells.Item(73, 13).Value = -239064  # M73
$ws.Cells.Item(73, 14).Value = -301872  # N73
$ws.Cells.Item(132, 8).Value = 1075.5135  # H132
$ws.Cells.Item(132, 9).Value = 1077.6111  # I132
$ws.Cells.Item(132, 11).Value = 3232.8333  # K132
$ws.Cells.Item(132, 13).Value = -702.8333000000002  # M132
$ws.Cells.Item(137, 8).Value = 1548.7368  # H137
$ws.Cells.Item(137, 9).Value = 1189.6364  # I137
$ws.Cells.Item(137, 11).Value = 3568.9092  # K137
$ws.Cells.Item(137, 13).Value = -1018.9092  # M137
$ws.Cells.Item(138, 8).Value = 2540  # H138
$ws.Cells.Item(138, 9).Value = 2842.25  # I138
$ws.Cells.Item(138, 10).Value = 2108.2144  # J138
$ws.Cells.Item(138, 11).Value = 8526.75  # K138
$ws.Cells.Item(138, 12).Value = 6324.6432  # L138
$ws.Cells.Item(138, 13).Value = -3386.75  # M138
$ws.Cells.Item(138, 14).Value = -16604.6432  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2935.5  # H32
$ws.Cells.Item(32, 9).Value = 1898.619  # I32
$ws.Cells.Item(32, 11).Value = 1898.619  # K32
$ws.Cells.Item(32, 13).Value = -1611.619  # M32
$ws.Cells.Item(74, 8).Value = 1793.375  # H74
$ws.Cells.Item(74, 10).Value = 5507  # J74
$ws.Cells.Item(74, 12).Value = 5507  # L74
$ws.Cells.Item(74, 14).Value = -7255  # N74
$ws.Cells.Item(77, 8).Value = 1793.375  # H77
$ws.Cells.Item(77, 10).Value = 5507  # J77
$ws.Cells.Item(77, 12).Value = 27535  # L77
$ws.Cells.Item(77, 14).Value = -36271  # N77
$ws.Cells.Item(110, 8).Value = 1481.0869  # H110
$ws.Cells.Item(110, 9).Value = 1187.2858  # I110
$ws.Cells.Item(110, 10).Value = 1938.1111  # J110
$ws.Cells.Item(110, 11).Value = 1187.2858  # K110
$ws.Cells.Item(110, 12).Value = 1938.1111  # L110
$ws.Cells.Item(110, 13).Value = 857.7141999999999  # M110
$ws.Cells.Item(110, 14).Value = -6028.1111  # N110

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value = 25000  # H100
$ws.Cells.Item(100, 10).Value = 25000  # J100
$ws.Cells.Item(100, 12).Value = 25000  # L100
$ws.Cells.Item(100, 14).Value = -27164  # N100
$ws.Cells.Item(105, 8).Value = 2506.6667  # H105
$ws.Cells.Item(105, 9).Value = 2506.6667  # I105
$ws.Cells.Item(105, 11).Value = 2506.6667  # K105
$ws.Cells.Item(105, 13).Value = -759.6667000000002  # M105
$ws.Cells.Item(107, 8).Value = 2395.077  # H107
$ws.Cells.Item(107, 9).Value = 2016.6666  # I107
$ws.Cells.Item(107, 11).Value = 2016.6666  # K107
$ws.Cells.Item(107, 13).Value = -96.66660000000002  # M107
$ws.Cells.Item(134, 8).Value = 7326.778  # H134
$ws.Cells.Item(134, 9).Value = 8460  # I134
$ws.Cells.Item(134, 11).Value = 25380  # K134
$ws.Cells.Item(134, 13).Value = -22845  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 52.57143  # H7
$ws.Cells.Item(7, 9).Value = 68.5  # I7
$ws.Cells.Item(7, 11).Value = 68.5  # K7
$ws.Cells.Item(7, 13).Value = 44.5  # M7
$ws.Cells.Item(22, 8).Value = 649.8  # H22
$ws.Cells.Item(22, 9).Value = 312.25  # I22
$ws.Cells.Item(22, 10).Value = 2000  # J22
$ws.Cells.Item(22, 11).Value = 312.25  # K22
$ws.Cells.Item(22, 12).Value = 2000  # L22
$ws.Cells.Item(22, 13).Value = 37.75  # M22
$ws.Cells.Item(22, 14).Value = -2700  # N22
$ws.Cells.Item(29, 8).Value = 0  # H29
$ws.Cells.Item(29, 10).Value = 0  # J29
$ws.Cells.Item(29, 14).ClearContents()  # N29
$ws.Cells.Item(31, 8).Value = 2728  # H31
$ws.Cells.Item(31, 9).Value = 1262.9  # I31
$ws.Cells.Item(31, 11).Value = 1262.9  # K31
$ws.Cells.Item(31, 13).Value = -967.9000000000001  # M31
$ws.Cells.Item(34, 8).Value = 2728  # H34
$ws.Cells.Item(34, 9).Value = 1262.9  # I34
$ws.Cells.Item(34, 11).Value = 1262.9  # K34
$ws.Cells.Item(34, 13).Value = -1060.9  # M34
$ws.Cells.Item(62, 8).Value = 5000  # H62
$ws.Cells.Item(62, 9).Value = 5000  # I62
$ws.Cells.Item(62, 11).Value = 5000  # K62
$ws.Cells.Item(62, 13).Value = -4376  # M62
$ws.Cells.Item(65, 8).Value = 5000  # H65
$ws.Cells.Item(65, 9).Value = 5000  # I65
$ws.Cells.Item(65, 11).Value = 25000  # K65
$ws.Cells.Item(65, 13).Value = -21880  # M65
$ws.Cells.Item(86, 8).Value = 1835.6666  # H86
$ws.Cells.Item(86, 9).Value = 1835.6666  # I86
$ws.Cells.Item(86, 11).Value = 1835.6666  # K86
$ws.Cells.Item(86, 13).Value = -712.6666  # M86
$ws.Cells.Item(89, 8).Value = 1835.6666  # H89
$ws.Cells.Item(89, 9).Value = 1835.6666  # I89
$ws.Cells.Item(89, 11).Value = 9178.333000000001  # K89
$ws.Cells.Item(89, 13).Value = -3562.333000000001  # M89
$ws.Cells.Item(107, 8).Value = 379.5  # H107
$ws.Cells.Item(107, 9).Value = 299.05  # I107
$ws.Cells.Item(107, 10).Value = 647.6667  # J107
$ws.Cells.Item(107, 11).Value = 299.05  # K107
$ws.Cells.Item(107, 12).Value = 647.6667  # L107
$ws.Cells.Item(107, 13).Value = 1620.95  # M107
$ws.Cells.Item(107, 14).Value = -4487.6667  # N107
$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 14).ClearContents()  # N130

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 238.6842  # H7
$ws.Cells.Item(7, 9).Value = 75.666664  # I7
$ws.Cells.Item(7, 11).Value = 226.999992  # K7
$ws.Cells.Item(7, 13).Value = -114.999992  # M7
$ws.Cells.Item(121, 8).Value = 529.2  # H121
$ws.Cells.Item(121, 10).Value = 750  # J121
$ws.Cells.Item(121, 12).Value = 2250  # L121
$ws.Cells.Item(121, 14).Value = -4870  # N121
$ws.Cells.Item(122, 8).Value = 871.3  # H122
$ws.Cells.Item(122, 10).Value = 1246.5  # J122
$ws.Cells.Item(122, 12).Value = 11218.5  # L122
$ws.Cells.Item(122, 14).Value = -16118.5  # N122
$ws.Cells.Item(131, 8).Value = 5961456  # H131
$ws.Cells.Item(131, 9).Value = 71429144  # I131
$ws.Cells.Item(131, 10).Value = 9848.169  # J131
$ws.Cells.Item(131, 11).Value = 214287432  # K131
$ws.Cells.Item(131, 12).Value = 29544.507  # L131
$ws.Cells.Item(131, 13).Value = -214282392  # M131
$ws.Cells.Item(131, 14).Value = -39624.507  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 85.833336  # H2
$ws.Cells.Item(2, 9).Value = 38.333332  # I2
$ws.Cells.Item(2, 10).Value = 133.33333  # J2
$ws.Cells.Item(2, 11).Value = 38.333332  # K2
$ws.Cells.Item(2, 12).Value = 133.33333  # L2
$ws.Cells.Item(2, 13).Value = 74.666668  # M2
$ws.Cells.Item(2, 14).Value = -359.33333  # N2
$ws.Cells.Item(102, 8).Value = 2764.923  # H102
$ws.Cells.Item(102, 9).Value = 3475  # I102
$ws.Cells.Item(102, 10).Value = 2156.2856  # J102
$ws.Cells.Item(102, 11).Value = 3475  # K102
$ws.Cells.Item(102, 12).Value = 2156.2856  # L102
$ws.Cells.Item(102, 13).Value = -1853  # M102
$ws.Cells.Item(102, 14).Value = -5400.2856  # N102
$ws.Cells.Item(122, 8).Value = 2205.4707  # H122
$ws.Cells.Item(122, 9).Value = 2027.7273  # I122
$ws.Cells.Item(122, 11).Value = 6083.1819  # K122
$ws.Cells.Item(122, 13).Value = -3633.1819  # M122
$ws.Cells.Item(123, 8).Value = 15195.4  # H123
$ws.Cells.Item(123, 10).Value = 15195.4  # J123
$ws.Cells.Item(123, 12).Value = 15195.4  # L123
$ws.Cells.Item(123, 14).Value = -20095.4  # N123
$ws.Cells.Item(126, 8).Value = 48468.316  # H126
$ws.Cells.Item(126, 9).Value = 3363.25  # I126
$ws.Cells.Item(126, 10).Value = 168748.5  # J126
$ws.Cells.Item(126, 11).Value = 10089.75  # K126
$ws.Cells.Item(126, 12).Value = 506245.5  # L126
$ws.Cells.Item(126, 13).Value = -7619.75  # M126
$ws.Cells.Item(126, 14).Value = -511185.5  # N126
$ws.Cells.Item(127, 8).Value = 36649  # H127
$ws.Cells.Item(127, 10).Value = 36649  # J127
$ws.Cells.Item(127, 12).Value = 36649  # L127
$ws.Cells.Item(127, 14).Value = -46569  # N127
$ws.Cells.Item(132, 8).Value = 3371.6829  # H132
$ws.Cells.Item(132, 9).Value = 2868.1936  # I132
$ws.Cells.Item(132, 10).Value = 4932.5  # J132
$ws.Cells.Item(132, 11).Value = 8604.5808  # K132
$ws.Cells.Item(132, 12).Value = 14797.5  # L132
$ws.Cells.Item(132, 13).Value = -6074.5808  # M132
$ws.Cells.Item(132, 14).Value = -19857.5  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1706.0834  # H22
$ws.Cells.Item(22, 10).Value = 1929.75  # J22
$ws.Cells.Item(22, 12).Value = 1929.75  # L22
$ws.Cells.Item(22, 14).Value = -2519.75  # N22
$ws.Cells.Item(27, 8).Value = 1706.0834  # H27
$ws.Cells.Item(27, 10).Value = 1929.75  # J27
$ws.Cells.Item(27, 12).Value = 1929.75  # L27
$ws.Cells.Item(27, 14).Value = -2143.75  # N27
$ws.Cells.Item(40, 8).Value = 12398.111  # H40
$ws.Cells.Item(40, 9).Value = 7799  # I40
$ws.Cells.Item(40, 11).Value = 7799  # K40
$ws.Cells.Item(40, 13).Value = -7663  # M40
$ws.Cells.Item(122, 8).Value = 7404.909  # H122
$ws.Cells.Item(122, 9).Value = 6010.3  # I122
$ws.Cells.Item(122, 10).Value = 8567.083000000001  # J122
$ws.Cells.Item(122, 11).Value = 18030.9  # K122
$ws.Cells.Item(122, 12).Value = 25701.249  # L122
$ws.Cells.Item(122, 13).Value = -15580.9  # M122
$ws.Cells.Item(122, 14).Value = -30601.249  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 623.7778  # H107
$ws.Cells.Item(107, 10).Value = 910  # J107
$ws.Cells.Item(107, 12).Value = 2730  # L107
$ws.Cells.Item(107, 14).Value = -6570  # N107
$ws.Cells.Item(132, 8).Value = 2549.96  # H132
$ws.Cells.Item(132, 9).Value = 2254.2632  # I132
$ws.Cells.Item(132, 11).Value = 6762.7896  # K132
$ws.Cells.Item(132, 13).Value = -4232.7896  # M132
